# "Fruta / hortaliza, semanal" — weekly refresh of the Zapallo price sheet.
#
# A new weekly observation is inserted at row 172 (pushing the existing
# rows 172-230 down to 173-231), growing the used range from A1:R230 to
# A1:R231. The new row carries this week's reading; every other row keeps
# its original data, just shifted down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 172; everything below (old 172..230) shifts to 173..231.
$ws.Rows("172:172").Insert()

$ws.Range("A172").Value = 4
$ws.Range("B172").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C172").Value = "Los Lagos"
$ws.Range("D172").Value = 44559
$ws.Range("E172").Value = 10
$ws.Range("F172").Value = 100112045
$ws.Range("G172").Value = "Zapallo"
$ws.Range("H172").Value = "Paine"
$ws.Range("I172").Value = "1a nueva(o)"
$ws.Range("J172").Value = 400
$ws.Range("K172").Value = 450
$ws.Range("L172").Value = 500
$ws.Range("M172").Value = 475
$ws.Range("N172").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O172").Value = "Región de O'Higgins"
$ws.Range("P172").Value = 475
$ws.Range("Q172").Value = 1
$ws.Range("R172").Value = "Hortaliza"
